$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Npy"
$ws.Range("C2").Value = "Fap"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.672624
$ws.Range("H2").Value = 2.017872
$ws.Range("I2").Value = 0.5229110735646243
$ws.Range("J2").Value = 0.5229110735646243
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.144859
$ws.Range("N2").Value = 6.434577
$ws.Range("O2").Value = 0.01238161622251906
$ws.Range("P2").Value = 0.01238161622251906
$ws.Range("Q2").Value = 1.442683640016
$ws.Range("R2").Value = 12.984152760144
$ws.Range("S2").Value = 0.006474484231382612
$ws.Range("T2").Value = 0.006474484231382612

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Npy"
$ws.Range("C3").Value = "Fap"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.672624
$ws.Range("H3").Value = 2.017872
$ws.Range("I3").Value = 0.5229110735646243
$ws.Range("J3").Value = 0.5229110735646243
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 155.9453076666667
$ws.Range("N3").Value = 467.835923
$ws.Range("O3").Value = 0.9002246540361518
$ws.Range("P3").Value = 0.9002246540361518
$ws.Range("Q3").Value = 104.892556623984
$ws.Range("R3").Value = 944.033009615856
$ws.Range("S3").Value = 0.4707374402913866
$ws.Range("T3").Value = 0.4707374402913866

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Npy"
$ws.Range("C4").Value = "Fap"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.672624
$ws.Range("H4").Value = 2.017872
$ws.Range("I4").Value = 0.5229110735646243
$ws.Range("J4").Value = 0.5229110735646243
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.823368666666667
$ws.Range("N4").Value = 8.470105999999999
$ws.Range("O4").Value = 0.0162984453921456
$ws.Range("P4").Value = 0.0162984453921456
$ws.Range("Q4").Value = 1.899065526048
$ws.Range("R4").Value = 17.091589734432
$ws.Range("S4").Value = 0.008522637577441261
$ws.Range("T4").Value = 0.008522637577441259

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Npy"
$ws.Range("C5").Value = "Fap"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.672624
$ws.Range("H5").Value = 2.017872
$ws.Range("I5").Value = 0.5229110735646243
$ws.Range("J5").Value = 0.5229110735646243
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 12.315788
$ws.Range("N5").Value = 36.947364
$ws.Range("O5").Value = 0.07109528434918362
$ws.Range("P5").Value = 0.07109528434918362
$ws.Range("Q5").Value = 8.283894587712
$ws.Range("R5").Value = 74.555051289408
$ws.Range("S5").Value = 0.03717651146441384
$ws.Range("T5").Value = 0.03717651146441384

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Npy"
$ws.Range("C6").Value = "Fap"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.6136826666666667
$ws.Range("H6").Value = 1.841048
$ws.Range("I6").Value = 0.4770889264353757
$ws.Range("J6").Value = 0.4770889264353757
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.144859
$ws.Range("N6").Value = 6.434577
$ws.Range("O6").Value = 0.01238161622251906
$ws.Range("P6").Value = 0.01238161622251906
$ws.Range("Q6").Value = 1.316262790744
$ws.Range("R6").Value = 11.846365116696
$ws.Range("S6").Value = 0.005907131991136452
$ws.Range("T6").Value = 0.005907131991136452

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Npy"
$ws.Range("C7").Value = "Fap"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.6136826666666667
$ws.Range("H7").Value = 1.841048
$ws.Range("I7").Value = 0.4770889264353757
$ws.Range("J7").Value = 0.4770889264353757
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 155.9453076666667
$ws.Range("N7").Value = 467.835923
$ws.Range("O7").Value = 0.9002246540361518
$ws.Range("P7").Value = 0.9002246540361518
$ws.Range("Q7").Value = 95.70093226303378
$ws.Range("R7").Value = 861.3083903673041
$ws.Range("S7").Value = 0.4294872137447652
$ws.Range("T7").Value = 0.4294872137447652

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Npy"
$ws.Range("C8").Value = "Fap"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.6136826666666667
$ws.Range("H8").Value = 1.841048
$ws.Range("I8").Value = 0.4770889264353757
$ws.Range("J8").Value = 0.4770889264353757
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.823368666666667
$ws.Range("N8").Value = 8.470105999999999
$ws.Range("O8").Value = 0.0162984453921456
$ws.Range("P8").Value = 0.0162984453921456
$ws.Range("Q8").Value = 1.732652412343111
$ws.Range("R8").Value = 15.593871711088
$ws.Range("S8").Value = 0.007775807814704341
$ws.Range("T8").Value = 0.00777580781470434

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Npy"
$ws.Range("C9").Value = "Fap"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.6136826666666667
$ws.Range("H9").Value = 1.841048
$ws.Range("I9").Value = 0.4770889264353757
$ws.Range("J9").Value = 0.4770889264353757
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 12.315788
$ws.Range("N9").Value = 36.947364
$ws.Range("O9").Value = 0.07109528434918362
$ws.Range("P9").Value = 0.07109528434918362
$ws.Range("Q9").Value = 7.557985621941333
$ws.Range("R9").Value = 68.021870597472
$ws.Range("S9").Value = 0.03391877288476979
$ws.Range("T9").Value = 0.03391877288476979
